$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin names, links) - no numeric coercion risk.
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +3.46%  "
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("E21").Value = "  +12.73%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("E29").Value = "  +15.11%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -4.44%  "
$ws.Range("E46").Value = "  +8.32%  "
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -3.57%  "

# Price cells (column D) - force text storage so numeric-looking
# strings ("312.05", "1.01", ...) are not coerced to numbers/floats,
# matching the source workbook where these are inline strings.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "42.718.60"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.280.68"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.05"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "101.90"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.622"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.595"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "38.65"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0897"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.20"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.974"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.03"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.632.88"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.282.77"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "42.694.67"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.28"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000104"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.25"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "73.10"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "261.35"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.65"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.81"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "22.22"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "36.06"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "166.26"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0858"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.62"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.50"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0346"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.62"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.63"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "69.18"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.226"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "91.42"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.720.56"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "11.89"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "110.43"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "78.11"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "8.65"
$c.Style = "Normal"
